$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 543.05554
$ws.Range("I19").Value = 447.13333
$ws.Range("J19").Value = 611.5714
$ws.Range("K19").Value = 447.13333
$ws.Range("L19").Value = 611.5714
$ws.Range("M19").Value = -272.13333
$ws.Range("N19").Value = -961.5714

$ws.Range("H43").Value = 712.46155
$ws.Range("I43").Value = 880
$ws.Range("K43").Value = 880
$ws.Range("M43").Value = -811

$ws.Range("H53").Value = 290.75
$ws.Range("I53").Value = 152.55556
$ws.Range("J53").Value = 468.42856
$ws.Range("K53").Value = 152.55556
$ws.Range("L53").Value = 468.42856
$ws.Range("M53").Value = 484.44444
$ws.Range("N53").Value = -1742.42856

$ws.Range("H116").Value = 23063582
$ws.Range("I116").Value = 27675298
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 27675298
$ws.Range("L116").Value = 5000
$ws.Range("M116").Value = -27671856
$ws.Range("N116").Value = -11884

$ws.Range("H121").Value = 525.1429000000001
$ws.Range("J121").Value = 525.1429000000001
$ws.Range("L121").Value = 1575.4287
$ws.Range("N121").Value = -5069.4287

$ws.Range("H138").Value = 17861038
$ws.Range("I138").Value = 6727.2856
$ws.Range("J138").Value = 23812476
$ws.Range("K138").Value = 20181.8568
$ws.Range("L138").Value = 71437428
$ws.Range("M138").Value = -15041.8568
$ws.Range("N138").Value = -71447708

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2217.8628
$ws.Range("I32").Value = 1741.381
$ws.Range("J32").Value = 4441.4443
$ws.Range("K32").Value = 1741.381
$ws.Range("L32").Value = 4441.4443
$ws.Range("M32").Value = -1454.381
$ws.Range("N32").Value = -5015.4443

$ws.Range("H122").Value = 1753.125
$ws.Range("I122").Value = 1527.4286
$ws.Range("J122").Value = 3333
$ws.Range("K122").Value = 4582.2858
$ws.Range("L122").Value = 9999
$ws.Range("M122").Value = -2132.2858
$ws.Range("N122").Value = -14899

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2015.5652
$ws.Range("I107").Value = 1966.3158
$ws.Range("K107").Value = 1966.3158
$ws.Range("M107").Value = -46.31580000000008

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3149.15
$ws.Range("I132").Value = 2165.7273
$ws.Range("K132").Value = 6497.1819
$ws.Range("M132").Value = -3967.1819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1668.8154
$ws.Range("I68").Value = 1818.5186
$ws.Range("J68").Value = 1562.4474
$ws.Range("K68").Value = 5455.5558
$ws.Range("L68").Value = 4687.3422
$ws.Range("M68").Value = -4644.5558
$ws.Range("N68").Value = -6309.3422

$ws.Range("H71").Value = 1668.8154
$ws.Range("I71").Value = 1818.5186
$ws.Range("J71").Value = 1562.4474
$ws.Range("K71").Value = 16366.6674
$ws.Range("L71").Value = 14062.0266
$ws.Range("M71").Value = -12310.6674
$ws.Range("N71").Value = -22174.0266

$ws.Range("H131").Value = 2338.5117
$ws.Range("I131").Value = 378.57144
$ws.Range("J131").Value = 2512.1772
$ws.Range("K131").Value = 1135.71432
$ws.Range("L131").Value = 7536.5316
$ws.Range("M131").Value = 3904.28568
$ws.Range("N131").Value = -17616.5316

$ws.Range("H132").Value = 957.4286
$ws.Range("I132").Value = 584
$ws.Range("J132").Value = 1237.5
$ws.Range("K132").Value = 5256
$ws.Range("L132").Value = 11137.5
$ws.Range("M132").Value = -2726
$ws.Range("N132").Value = -16197.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 1201100
$ws.Range("J7").Value = 2750
$ws.Range("L7").Value = 2750
$ws.Range("N7").Value = -2974

$ws.Range("H8").Value = 1201100
$ws.Range("J8").Value = 2750
$ws.Range("L8").Value = 2750
$ws.Range("N8").Value = -3028

$ws.Range("H107").Value = 2020
$ws.Range("I107").Value = 2020
$ws.Range("K107").Value = 2020
$ws.Range("M107").Value = -100

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").Value = $null

$ws.Range("H113").Value = 1067.4286
$ws.Range("I113").Value = 745.3333
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 745.3333
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 1424.6667
$ws.Range("N113").Value = -7340

$ws.Range("H122").Value = 1591417
$ws.Range("J122").Value = 2504
$ws.Range("L122").Value = 7512
$ws.Range("N122").Value = -12412

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 60001.5
$ws.Range("I2").Value = 50001
$ws.Range("K2").Value = 50001
$ws.Range("M2").Value = -49889

$ws.Range("H7").Value = 3096.1853
$ws.Range("I7").Value = 1879.8
$ws.Range("K7").Value = 1879.8
$ws.Range("M7").Value = -1767.8

$ws.Range("H16").Value = 1430.3334
$ws.Range("I16").Value = 1545
$ws.Range("J16").Value = 1201
$ws.Range("K16").Value = 1545
$ws.Range("L16").Value = 1201
$ws.Range("M16").Value = -1375
$ws.Range("N16").Value = -1541

$ws.Range("H40").Value = 4458.4707
$ws.Range("I40").Value = 2931.3333
$ws.Range("J40").Value = 4785.7144
$ws.Range("K40").Value = 2931.3333
$ws.Range("L40").Value = 4785.7144
$ws.Range("M40").Value = -2795.3333
$ws.Range("N40").Value = -5057.7144

$ws.Range("H61").Value = 2088.625
$ws.Range("I61").Value = 2245.5715
$ws.Range("J61").Value = 990
$ws.Range("K61").Value = 2245.5715
$ws.Range("L61").Value = 990
$ws.Range("M61").Value = -2043.5715
$ws.Range("N61").Value = -1394

$ws.Range("H113").Value = 2088.625
$ws.Range("I113").Value = 2245.5715
$ws.Range("J113").Value = 990
$ws.Range("K113").Value = 2245.5715
$ws.Range("L113").Value = 990
$ws.Range("M113").Value = -75.57150000000001
$ws.Range("N113").Value = -5330

$ws.Range("H126").Value = 3096.1853
$ws.Range("I126").Value = 1879.8
$ws.Range("K126").Value = 5639.4
$ws.Range("M126").Value = -3169.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 11344.4
$ws.Range("J74").Value = 11038.25
$ws.Range("L74").Value = 11038.25
$ws.Range("N74").Value = -12910.25

$ws.Range("H77").Value = 11344.4
$ws.Range("J77").Value = 11038.25
$ws.Range("L77").Value = 33114.75
$ws.Range("N77").Value = -42474.75

$ws.Range("H132").Value = 13159320
$ws.Range("I132").Value = 18519780
$ws.Range("J132").Value = 1828.909
$ws.Range("K132").Value = 55559340
$ws.Range("L132").Value = 5486.727000000001
$ws.Range("M132").Value = -55556810
$ws.Range("N132").Value = -10546.727
